# Update "想去人数" (interest count) values in column F across the
# workbook's sheets, reflecting refreshed scrape numbers.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 220
$ws.Range("F3").Value = 54857
$ws.Range("F4").Value = 1347
$ws.Range("F6").Value = 349
$ws.Range("F7").Value = 885
$ws.Range("F8").Value = 776
$ws.Range("F9").Value = 415
$ws.Range("F10").Value = 3092
$ws.Range("F12").Value = 5256
$ws.Range("F14").Value = 1096
$ws.Range("F18").Value = 427
$ws.Range("F19").Value = 1320
$ws.Range("F22").Value = 194
$ws.Range("F24").Value = 39
$ws.Range("F29").Value = 5195
$ws.Range("F31").Value = 5121
$ws.Range("F32").Value = 9120
$ws.Range("F35").Value = 143
$ws.Range("F40").Value = 4236
$ws.Range("F41").Value = 260

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 101
$ws.Range("F5").Value = 135
$ws.Range("F12").Value = 1141
$ws.Range("F18").Value = 57

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 786
$ws.Range("F3").Value = 576

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 786
$ws.Range("F3").Value = 576
$ws.Range("F4").Value = 220
$ws.Range("F5").Value = 349
$ws.Range("F6").Value = 885
$ws.Range("F7").Value = 776
$ws.Range("F8").Value = 415
$ws.Range("F9").Value = 3093
$ws.Range("F11").Value = 101
$ws.Range("F13").Value = 135
$ws.Range("F14").Value = 1096
$ws.Range("F18").Value = 427
$ws.Range("F20").Value = 1320
$ws.Range("F23").Value = 194
$ws.Range("F28").Value = 5195
$ws.Range("F29").Value = 5121
$ws.Range("F30").Value = 9120
$ws.Range("F34").Value = 143
$ws.Range("F41").Value = 4236
$ws.Range("F42").Value = 57
$ws.Range("F47").Value = 260
